# 角色.xlsx — "Add files via upload" edit:
#  - 功能面板!B9  1 -> 0   (also recalcs paralimit!B3)
#  - 功能面板!B11 0.05 -> 0 (also recalcs xzqlimit!B2)
#  - selection / active-sheet bookkeeping changes:
#      功能面板 loses the "active tab" + selection moves to B12
#      基础面板 selection moves to E28
#      calc     becomes the active tab + selection moves to F18
#  - paralimit sheet picks up an explicit Page Setup (paper size 9 / portrait)

$wb = $excel.ActiveWorkbook

# --- 功能面板 (1st sheet): core value edits ---
$wsFunc = $wb.Worksheets.Item(1)
$wsFunc.Range("B9").Value = 0
$wsFunc.Range("B11").Value = 0

# --- 基础面板 (3rd sheet): selection moves ---
$wsBase = $wb.Worksheets.Item(3)
$wsBase.Activate() | Out-Null
$wsBase.Range("E28").Select() | Out-Null

# go back and park the selection on 功能面板 before handing off the
# "active tab" to calc (matches the final selection recorded for it)
$wsFunc.Activate() | Out-Null
$wsFunc.Range("B12").Select() | Out-Null

# --- paralimit (7th sheet): give it an explicit page setup ---
$wsPara = $wb.Worksheets.Item(7)
$wsPara.PageSetup.PaperSize = 9
$wsPara.PageSetup.Orientation = 1

# --- calc (6th sheet): ends up the active / tab-selected sheet ---
$wsCalc = $wb.Worksheets.Item(6)
$wsCalc.Activate() | Out-Null
$wsCalc.Range("F18").Select() | Out-Null
